# Refresh the "quadratic-svm-score" prediction-distance value (B2) with the
# latest figure from the re-run pipeline (previous copy of ful-path.csv).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quadratic-svm-score")
if (-not $ws) {
    $ws = $wb.ActiveSheet
}

$ws.Range("B2").Value = 1684.2214090878099
